# Implement basket-based elective scheduling with common time slots across all branches
# Updates the timetable cell values on both Section_A and Section_B sheets so that
# elective slots (CS307 / EC262 / EC264 / EC303 / EC304 and their tutorials) line up
# on common time slots across branches/sections.

$wb = $excel.ActiveWorkbook

$sectionA = $wb.Worksheets.Item("Section_A")
$sectionB = $wb.Worksheets.Item("Section_B")

$updatesA = @{
    "B2" = "EC303"
    "C2" = "Free"
    "D2" = "EC264"
    "E2" = "Free"
    "F2" = "EC304"
    "C3" = "Free"
    "E3" = "CS307"
    "F3" = "EC264"
    "B5" = "CS307"
    "C5" = "Free"
    "F5" = "EC262"
    "D6" = "CS307 (Tutorial)"
    "F6" = "Free"
    "B7" = "EC304"
    "C7" = "CS307"
    "D7" = "Free"
    "E7" = "EC303"
    "F7" = "Free"
    "D8" = "Free"
    "F8" = "Free"
}

$updatesB = @{
    "B2" = "CS307"
    "C2" = "Free"
    "D2" = "CS307"
    "E2" = "Free"
    "F2" = "EC262"
    "B3" = "EC304"
    "C3" = "EC303"
    "D3" = "Free"
    "F3" = "EC303"
    "B5" = "Free"
    "C5" = "Free"
    "D5" = "EC262"
    "E5" = "EC304"
    "F5" = "EC264"
    "C6" = "Free"
    "E6" = "CS307 (Tutorial)"
    "F6" = "Free"
    "B7" = "EC264"
    "C7" = "Free"
    "D7" = "EC304"
    "E7" = "CS307"
    "F7" = "Free"
    "C8" = "Free"
}

foreach ($addr in $updatesA.Keys) {
    $sectionA.Range($addr).Value = $updatesA[$addr]
}

foreach ($addr in $updatesB.Keys) {
    $sectionB.Range($addr).Value = $updatesB[$addr]
}
